$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "27.680.13"
Set-TextValue "E2" "  +0.00%  "
Set-TextValue "D3" "1.584.14"
Set-TextValue "E3" "  -2.31%  "
Set-TextValue "E4" "  +1.05%  "
Set-TextValue "D5" "206.90"
Set-TextValue "E5" "  -1.47%  "
Set-TextValue "D6" "0.506"
Set-TextValue "E6" "  -1.97%  "
Set-TextValue "E7" "  +1.05%  "
Set-TextValue "D8" "22.25"
Set-TextValue "E8" "  -4.09%  "
Set-TextValue "E9" "  -1.32%  "
Set-TextValue "D10" "0.0590"
Set-TextValue "E10" "  -2.67%  "
Set-TextValue "E11" "  -0.92%  "
Set-TextValue "D12" "1.809.30"
Set-TextValue "E12" "  -2.31%  "
Set-TextValue "D13" "1.585.30"
Set-TextValue "E13" "  -3.43%  "
Set-TextValue "E14" "  -3.12%  "
Set-TextValue "D15" "0.530"
Set-TextValue "E15" "  -5.01%  "
Set-TextValue "D16" "27.657.57"
Set-TextValue "E16" "  -0.19%  "
Set-TextValue "D17" "63.22"
Set-TextValue "E17" "  -2.41%  "
Set-TextValue "D18" "219.00"
Set-TextValue "E18" "  -3.77%  "
Set-TextValue "D19" "0.0₃0692"
Set-TextValue "E19" "  -3.22%  "
Set-TextValue "E20" "  -4.60%  "
Set-TextValue "E21" "  +0.93%  "
Set-TextValue "E22" "  -4.11%  "
Set-TextValue "E23" "  -5.78%  "
Set-TextValue "E24" "  -3.76%  "
Set-TextValue "D25" "154.59"
Set-TextValue "E25" "  -0.01%  "
Set-TextValue "D26" "6.81"
Set-TextValue "E26" "  -1.48%  "
Set-TextValue "E27" "  +1.03%  "
Set-TextValue "D28" "15.11"
Set-TextValue "E28" "  -2.18%  "
Set-TextValue "E29" "  -3.35%  "
Set-TextValue "D30" "1.16"
Set-TextValue "E30" "  -1.45%  "
Set-TextValue "D31" "0.0464"
Set-TextValue "E31" "  -2.76%  "
Set-TextValue "E32" "  -4.72%  "
Set-TextValue "D33" "1.382.53"
Set-TextValue "E33" "  -0.70%  "
Set-TextValue "E34" "  -5.21%  "
Set-TextValue "E35" "  -4.41%  "
Set-TextValue "E36" "  -4.05%  "
Set-TextValue "E37" "  -0.16%  "
Set-TextValue "E38" "  -2.83%  "
Set-TextValue "D39" "0.537"
Set-TextValue "E39" "  -3.54%  "
Set-TextValue "D40" "0.820"
Set-TextValue "E40" "  -3.04%  "
Set-TextValue "E41" "  +0.89%  "
Set-TextValue "D42" "0.973"
Set-TextValue "E42" "  -3.84%  "
Set-TextValue "B43" "MXToken"
Set-TextValue "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.17"
Set-TextValue "E43" "  +0.44%  "
Set-TextValue "B44" "Aave"
Set-TextValue "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "63.50"
Set-TextValue "E44" "  -3.19%  "
Set-TextValue "E45" "  -4.15%  "
Set-TextValue "D46" "5.21"
Set-TextValue "E46" "  -3.39%  "
Set-TextValue "D47" "1.720.29"
Set-TextValue "E47" "  -2.63%  "
Set-TextValue "D48" "88.10"
Set-TextValue "E48" "  +0.26%  "
Set-TextValue "E49" "  +4.59%  "
Set-TextValue "D50" "0.0971"
Set-TextValue "E50" "  -4.34%  "
Set-TextValue "D51" "0.0499"
Set-TextValue "E51" "  -0.72%  "
